# Apply the cell-value updates to Sheet1 as described by the diff.
# Each row in the crypto table gets refreshed Price (D) / Volume(1h) (E)
# figures; rows 36-37 additionally swap which coin (Huobi/Lido) occupies
# which row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.543.92'
$ws.Range("E2").Value = '  -2.42%  '

$ws.Range("D3").Value = '1.582.93'
$ws.Range("E3").Value = '  -3.03%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.80'
$ws.Range("E5").Value = '  -2.73%  '

$ws.Range("E6").Value = '  -2.27%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  -2.79%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0618'

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0833'
$ws.Range("E11").Value = '  -1.77%  '

$ws.Range("D12").Value = '1.804.62'
$ws.Range("E12").Value = '  -3.03%  '

$ws.Range("D13").Value = '1.587.88'
$ws.Range("E13").Value = '  -2.85%  '

$ws.Range("E14").Value = '  -1.90%  '

$ws.Range("E15").Value = '  -2.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.36'

$ws.Range("D17").Value = '26.578.77'
$ws.Range("E17").Value = '  -2.19%  '

$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.34'
$ws.Range("E19").Value = '  -3.29%  '

$ws.Range("E20").Value = '  -0.06%  '

$ws.Range("E21").Value = '  -2.86%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.24'
$ws.Range("E22").Value = '  -3.56%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.39'
$ws.Range("E23").Value = '  -4.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.88'
$ws.Range("E24").Value = '  -2.34%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.48'
$ws.Range("E25").Value = '  -1.19%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.42'
$ws.Range("E26").Value = '  +1.84%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("E28").Value = '  -4.80%  '

$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("E30").Value = '  -0.56%  '

$ws.Range("E31").Value = '  -2.29%  '

$ws.Range("E32").Value = '  -3.87%  '

$ws.Range("E33").Value = '  +23.94%  '

$ws.Range("E34").Value = '  -2.65%  '

$ws.Range("D35").Value = '1.320.08'
$ws.Range("E35").Value = '  +0.61%  '

$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.43'
$ws.Range("E36").Value = '  -0.70%  '

$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.50'
$ws.Range("E37").Value = '  -4.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0172'
$ws.Range("E38").Value = '  -1.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.822'
$ws.Range("E39").Value = '  -3.45%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("E41").Value = '  -2.59%  '

$ws.Range("E42").Value = '  -3.48%  '

$ws.Range("E43").Value = '  +1.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.04'
$ws.Range("E44").Value = '  -1.39%  '

$ws.Range("D45").Value = '1.718.32'
$ws.Range("E45").Value = '  -2.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.04'
$ws.Range("E46").Value = '  -1.89%  '

$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.831'
$ws.Range("E48").Value = '  +2.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0988'
$ws.Range("E49").Value = '  +4.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("E50").Value = '  -1.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.47'
$ws.Range("E51").Value = '  -1.09%  '
